# Update AgTests (H) and AgPosit (I) columns for the Covid daily stats sheet
# reflecting the "ut 29. 12. 2020" data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H272").Value = 30655

$ws.Range("H273").Value = 26904
$ws.Range("I273").Value = 1366

$ws.Range("H274").Value = 28359
$ws.Range("I274").Value = 1338

$ws.Range("H275").Value = 28773
$ws.Range("I275").Value = 1238

$ws.Range("H278").Value = 29664
$ws.Range("I278").Value = 2056

$ws.Range("H279").Value = 42419
$ws.Range("I279").Value = 3086

$ws.Range("H280").Value = 36044
$ws.Range("I280").Value = 2397

$ws.Range("H281").Value = 45260
$ws.Range("I281").Value = 3262

$ws.Range("H282").Value = 46722
$ws.Range("I282").Value = 2841

$ws.Range("H286").Value = 54599
$ws.Range("I286").Value = 4223

$ws.Range("H287").Value = 56340
$ws.Range("I287").Value = 3842

$ws.Range("H288").Value = 55543
$ws.Range("I288").Value = 3939

$ws.Range("H289").Value = 63085
$ws.Range("I289").Value = 3575

$ws.Range("H290").Value = 17655
$ws.Range("I290").Value = 1480

$ws.Range("H292").Value = 77798
$ws.Range("I292").Value = 6926

$ws.Range("H293").Value = 78690
$ws.Range("I293").Value = 5669

$ws.Range("H294").Value = 87303
$ws.Range("I294").Value = 4828

$ws.Range("H295").Value = 19248

$ws.Range("H296").Value = 1863

$ws.Range("H297").Value = 2572

$ws.Range("H298").Value = 2564
$ws.Range("I298").Value = 268

$ws.Range("H299").Value = 57834
$ws.Range("I299").Value = 6074

$wb.Save()
